# Generate Report for Handback
# Updates the timestamp values that were refreshed when the handback report
# was regenerated.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-24 09:50:52"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-10-24 09:50:41"
$wsZhCn.Range("K2").Value = "2016-10-24 09:51:28"

# "de-de" sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-10-24 09:50:52"
$wsDeDe.Range("K2").Value = "2016-10-24 09:51:45"
